# Update cryptos list cells per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.204.68'
$ws.Range("E2").Value = '  -3.33%  '

$ws.Range("D3").Value = '2.225.26'
$ws.Range("E3").Value = '  -4.83%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.52'
$ws.Range("E5").Value = '  -2.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.28'
$ws.Range("E6").Value = '  -6.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.581'
$ws.Range("E7").Value = '  -7.67%  '

$ws.Range("E8").Value = '  -0.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  -7.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.06'
$ws.Range("E10").Value = '  -7.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.08'
$ws.Range("E11").Value = '  -2.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0828'
$ws.Range("E12").Value = '  -9.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.70'
$ws.Range("E13").Value = '  -7.50%  '

$ws.Range("E14").Value = '  -1.68%  '

$ws.Range("D15").Value = '2.561.91'
$ws.Range("E15").Value = '  -5.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.861'
$ws.Range("E16").Value = '  -10.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.33'
$ws.Range("E17").Value = '  -5.40%  '

$ws.Range("D18").Value = '2.220.79'
$ws.Range("E18").Value = '  -4.84%  '

$ws.Range("D19").Value = '43.084.76'
$ws.Range("E19").Value = '  -3.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.79'
$ws.Range("E20").Value = '  -9.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.57'
$ws.Range("E21").Value = '  -8.94%  '

$ws.Range("D22").Value = '0.0₃0964'
$ws.Range("E22").Value = '  -8.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.24'
$ws.Range("E23").Value = '  -10.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.14'
$ws.Range("E24").Value = '  -10.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '236.51'
$ws.Range("E25").Value = '  -7.52%  '

$ws.Range("E26").Value = '  -3.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.05'
$ws.Range("E28").Value = '  +1.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.02'
$ws.Range("E29").Value = '  -10.05%  '

$ws.Range("E30").Value = '  -2.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.43'
$ws.Range("E31").Value = '  -12.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.83'
$ws.Range("E32").Value = '  +1.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.23'
$ws.Range("E33").Value = '  -7.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0864'
$ws.Range("E34").Value = '  -8.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '157.16'
$ws.Range("E35").Value = '  -5.21%  '

$ws.Range("E36").Value = '  +3.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.68'
$ws.Range("E37").Value = '  -3.36%  '

$ws.Range("E38").Value = '  -7.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.85'
$ws.Range("E39").Value = '  -2.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.41'
$ws.Range("E40").Value = '  -5.37%  '

$ws.Range("E41").Value = '  -9.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.69'
$ws.Range("E42").Value = '  -5.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0318'
$ws.Range("E43").Value = '  -8.60%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.94'
$ws.Range("E44").Value = '  +9.71%  '

$ws.Range("E45").Value = '  -0.31%  '

$ws.Range("D46").Value = '1.751.93'
$ws.Range("E46").Value = '  -6.24%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.202'
$ws.Range("E47").Value = '  -9.67%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.95'
$ws.Range("E48").Value = '  -3.07%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.46'
$ws.Range("E49").Value = '  -12.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.26'
$ws.Range("E50").Value = '  -12.33%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.86'
$ws.Range("E51").Value = '  -11.08%  '
